# Scheduled-runner style update of price/profit figures across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1103.75
$ws.Range("I6").Value = 1304.5
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 3913.5
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -3801.5

$ws.Range("H15").Value = 1858.1666
$ws.Range("I15").Value = 1858.1666
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5574.4998
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -5405.4998

$ws.Range("H17").Value = 880086.6
$ws.Range("I17").Value = 1040
$ws.Range("J17").Value = 1190338.4
$ws.Range("K17").Value = 3120
$ws.Range("L17").Value = 3571015.2
$ws.Range("M17").Value = -2952
$ws.Range("N17").Value = -3571351.2

$ws.Range("H112").Value = 2495433.5
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2495433.5
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 7486300.5
$ws.Range("N112").Value = -7488516.5
$ws.Range("M112").ClearContents()

$ws.Range("H137").Value = 28685.844
$ws.Range("I137").Value = 32339
$ws.Range("J137").Value = 3113.75
$ws.Range("K137").Value = 97017
$ws.Range("L137").Value = 9341.25
$ws.Range("M137").Value = -94467
$ws.Range("N137").Value = -14441.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16198140
$ws.Range("I32").Value = 14468495
$ws.Range("J32").Value = 35718420
$ws.Range("K32").Value = 14468495
$ws.Range("L32").Value = 35718420
$ws.Range("M32").Value = -14468208
$ws.Range("N32").Value = -35718994

$ws.Range("H61").Value = 2442.4395
$ws.Range("I61").Value = 2234.9832
$ws.Range("J61").Value = 4191
$ws.Range("K61").Value = 2234.9832
$ws.Range("L61").Value = 4191
$ws.Range("M61").Value = -2022.9832

$ws.Range("H122").Value = 3935.4358
$ws.Range("I122").Value = 2159.32
$ws.Range("J122").Value = 7107.0713
$ws.Range("K122").Value = 6477.960000000001
$ws.Range("L122").Value = 21321.2139
$ws.Range("M122").Value = -4027.960000000001
$ws.Range("N122").Value = -26221.2139

$ws.Range("H132").Value = 105084.65
$ws.Range("I132").Value = 128214.516
$ws.Range("J132").Value = 3570.2222
$ws.Range("K132").Value = 384643.548
$ws.Range("L132").Value = 10710.6666
$ws.Range("M132").Value = -382113.548
$ws.Range("N132").Value = -15770.6666

$ws.Range("H136").Value = 2442.4395
$ws.Range("I136").Value = 2234.9832
$ws.Range("J136").Value = 4191
$ws.Range("K136").Value = 6704.9496
$ws.Range("L136").Value = 12573
$ws.Range("M136").Value = -4154.9496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1978.8518
$ws.Range("I86").Value = 1791.1333
$ws.Range("J86").Value = 2213.5
$ws.Range("K86").Value = 1791.1333
$ws.Range("L86").Value = 2213.5
$ws.Range("M86").Value = -668.1333
$ws.Range("N86").Value = -4459.5

$ws.Range("H89").Value = 1978.8518
$ws.Range("I89").Value = 1791.1333
$ws.Range("J89").Value = 2213.5
$ws.Range("K89").Value = 8955.666499999999
$ws.Range("L89").Value = 11067.5
$ws.Range("M89").Value = -3339.666499999999
$ws.Range("N89").Value = -22299.5

$ws.Range("H117").Value = 116940
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 116940
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 116940
$ws.Range("N117").Value = -126118

$ws.Range("H134").Value = 2103692.5
$ws.Range("I134").Value = 2553769.5
$ws.Range("J134").Value = 3333.1667
$ws.Range("K134").Value = 7661308.5
$ws.Range("L134").Value = 9999.500100000001
$ws.Range("M134").Value = -7658773.5
$ws.Range("N134").Value = -15069.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3550.3333
$ws.Range("I31").Value = 2818.6667
$ws.Range("J31").Value = 3916.1667
$ws.Range("K31").Value = 2818.6667
$ws.Range("L31").Value = 3916.1667
$ws.Range("M31").Value = -2523.6667
$ws.Range("N31").Value = -4506.1667

$ws.Range("H34").Value = 3550.3333
$ws.Range("I34").Value = 2818.6667
$ws.Range("J34").Value = 3916.1667
$ws.Range("K34").Value = 2818.6667
$ws.Range("L34").Value = 3916.1667
$ws.Range("M34").Value = -2616.6667
$ws.Range("N34").Value = -4320.1667

$ws.Range("H58").Value = 2308.2754
$ws.Range("I58").Value = 2085.4
$ws.Range("J58").Value = 2894.7896
$ws.Range("K58").Value = 2085.4
$ws.Range("L58").Value = 2894.7896
$ws.Range("M58").Value = -1882.4
$ws.Range("N58").Value = -3300.7896

$ws.Range("H99").Value = 3720.5
$ws.Range("I99").Value = 3542.5
$ws.Range("J99").Value = 3987.5
$ws.Range("K99").Value = 3542.5
$ws.Range("L99").Value = 3987.5
$ws.Range("M99").Value = -2044.5
$ws.Range("N99").Value = -6983.5

$ws.Range("H107").Value = 46924.953
$ws.Range("I107").Value = 72309.78999999999
$ws.Range("J107").Value = 2501.5
$ws.Range("K107").Value = 72309.78999999999
$ws.Range("L107").Value = 2501.5
$ws.Range("M107").Value = -70389.78999999999
$ws.Range("N107").Value = -6341.5

$ws.Range("H126").Value = 3720.5
$ws.Range("I126").Value = 3542.5
$ws.Range("J126").Value = 3987.5
$ws.Range("K126").Value = 10627.5
$ws.Range("L126").Value = 11962.5
$ws.Range("M126").Value = -8157.5
$ws.Range("N126").Value = -16902.5

$ws.Range("H132").Value = 3719.5942
$ws.Range("I132").Value = 3417.7166
$ws.Range("J132").Value = 5732.1113
$ws.Range("K132").Value = 10253.1498
$ws.Range("L132").Value = 17196.3339
$ws.Range("M132").Value = -7723.149800000001
$ws.Range("N132").Value = -22256.3339

$ws.Range("H136").Value = 2308.2754
$ws.Range("I136").Value = 2085.4
$ws.Range("J136").Value = 2894.7896
$ws.Range("K136").Value = 6256.200000000001
$ws.Range("L136").Value = 8684.3688
$ws.Range("M136").Value = -3706.200000000001
$ws.Range("N136").Value = -13784.3688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 671
$ws.Range("I5").Value = 516
$ws.Range("J5").Value = 1187.6666
$ws.Range("K5").Value = 1548
$ws.Range("L5").Value = 3562.9998
$ws.Range("M5").Value = -1436
$ws.Range("N5").Value = -3786.9998

$ws.Range("H9").Value = 1000.5
$ws.Range("I9").Value = 1999
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 5997
$ws.Range("L9").Value = 6
$ws.Range("M9").Value = -5773
$ws.Range("N9").Value = -454

$ws.Range("H131").Value = 1768.3823
$ws.Range("I131").Value = 1268.2858
$ws.Range("J131").Value = 1898.037
$ws.Range("K131").Value = 3804.8574
$ws.Range("L131").Value = 5694.111
$ws.Range("M131").Value = 1235.1426
$ws.Range("N131").Value = -15774.111

$ws.Range("H132").Value = 1283.0834
$ws.Range("I132").Value = 1320
$ws.Range("J132").Value = 1256.7142
$ws.Range("K132").Value = 11880
$ws.Range("L132").Value = 11310.4278
$ws.Range("M132").Value = -9350
$ws.Range("N132").Value = -16370.4278

$ws.Range("H135").Value = 671
$ws.Range("I135").Value = 516
$ws.Range("J135").Value = 1187.6666
$ws.Range("K135").Value = 4644
$ws.Range("L135").Value = 10688.9994
$ws.Range("M135").Value = -2109
$ws.Range("N135").Value = -15758.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3436.7144
$ws.Range("I126").Value = 4011.75
$ws.Range("J126").Value = 2670
$ws.Range("K126").Value = 12035.25
$ws.Range("L126").Value = 8010
$ws.Range("M126").Value = -9565.25
$ws.Range("N126").Value = -12950

$ws.Range("H132").Value = 2630
$ws.Range("I132").Value = 3384.4
$ws.Range("J132").Value = 2210.889
$ws.Range("K132").Value = 10153.2
$ws.Range("L132").Value = 6632.667
$ws.Range("M132").Value = -7623.200000000001
$ws.Range("N132").Value = -11692.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4022.1333
$ws.Range("I100").Value = 3985.6667
$ws.Range("J100").Value = 4168
$ws.Range("K100").Value = 3985.6667
$ws.Range("L100").Value = 4168
$ws.Range("M100").Value = -3444.6667
$ws.Range("N100").Value = -5250

$ws.Range("H132").Value = 837063.7
$ws.Range("I132").Value = 2502948.5
$ws.Range("J132").Value = 4121.25
$ws.Range("K132").Value = 7508845.5
$ws.Range("L132").Value = 12363.75
$ws.Range("M132").Value = -7506315.5
$ws.Range("N132").Value = -17423.75

$ws.Range("H136").Value = 2242.6538
$ws.Range("I136").Value = 1467.0952
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 4401.2856
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -1851.2856
$ws.Range("N136").Value = -21600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5294.091
$ws.Range("I62").Value = 3933
$ws.Range("J62").Value = 6071.857
$ws.Range("K62").Value = 3933
$ws.Range("L62").Value = 6071.857
$ws.Range("M62").Value = -3309
$ws.Range("N62").Value = -7319.857

$ws.Range("H65").Value = 5294.091
$ws.Range("I65").Value = 3933
$ws.Range("J65").Value = 6071.857
$ws.Range("K65").Value = 19665
$ws.Range("L65").Value = 30359.285
$ws.Range("M65").Value = -16545
$ws.Range("N65").Value = -36599.285

$ws.Range("H113").Value = 31251054
$ws.Range("I113").Value = 1155.3334
$ws.Range("J113").Value = 125000750
$ws.Range("K113").Value = 3466.0002
$ws.Range("L113").Value = 375002250
$ws.Range("M113").Value = -1296.0002
$ws.Range("N113").Value = -375006590

$ws.Range("H124").Value = 122370.664
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 122370.664
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 122370.664
$ws.Range("N124").Value = -132190.664

$ws.Range("H126").Value = 5206.2085
$ws.Range("I126").Value = 5543.1816
$ws.Range("J126").Value = 1499.5
$ws.Range("K126").Value = 16629.5448
$ws.Range("L126").Value = 4498.5
$ws.Range("M126").Value = -14159.5448
$ws.Range("N126").Value = -9438.5

$ws.Range("H132").Value = 31135.389
$ws.Range("I132").Value = 43265.16
$ws.Range("J132").Value = 3567.7273
$ws.Range("K132").Value = 129795.48
$ws.Range("L132").Value = 10703.1819
$ws.Range("M132").Value = -127265.48
$ws.Range("N132").Value = -15763.1819
